$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, mirroring the style of the existing G1 header cell.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells H2/H3 with plain numeric values (no special style).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
